# daily auto push: 2025-10-06 07:28 UTC
# Append the day's new log row (2025/10/06, 月, 16:00, ranking 201) as row 70
# right after the existing last row (69) of Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 70

# Column A holds a date-like literal ("2025/10/06") that must stay plain text
# (matching every other row in the sheet) instead of being auto-converted to
# a date serial number. Format the cell as Text before typing the value, then
# restore the default "Normal" style so no stray formatting is introduced.
$cellA = $ws.Cells.Item($row, 1)
$cellA.NumberFormat = "@"
$cellA.Value = "2025/10/06"
$cellA.NumberFormat = "General"
$cellA.Style = "Normal"

$ws.Cells.Item($row, 2).Value = "月"
$ws.Cells.Item($row, 3).Value = 16
$ws.Cells.Item($row, 4).Value = 201
